# Adapt tests to control version
# Adds a "version" column (header + value 1) to the "settings" sheet and
# makes that sheet the active/selected one (mirrors the author switching
# focus to the settings tab after adding the new control-version field).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# New third column: header "version" in row 1, value 1 in row 2.
$ws.Cells.Item(1, 3).Value = "version"
$ws.Cells.Item(2, 3).Value = 1

# Make "settings" the active sheet/tab and park the selection on the next
# empty cell (C3), matching where the cursor lands after filling C1:C2.
$ws.Activate()
$ws.Range("C3").Select()
